# Anu - Cash Management Files Uploaded
# Remove the sensitive credential data (URL / UserName / Password sample
# values) that were previously stored on the Input_Value sheet, row 2,
# columns N:P, along with the hyperlink that was attached to the URL cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input_Value")
$ws.Activate()

# The URL cell (N2) carries a hyperlink to the Oracle Cloud URL - drop it
# before clearing the cell text itself.
if ($ws.Range("N2").Hyperlinks.Count -gt 0) {
    $ws.Range("N2").Hyperlinks.Delete()
}

# Clear the credential values (URL, UserName, Password) out of N2:P2.
$ws.Range("N2:P2").ClearContents()

# Reflect the cells that were touched in the saved selection state.
$ws.Range("N2:P2").Select()
